$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Faculty" header in column G (new shared string, new cell)
$ws.Range("G1").Value = "Faculty"

# Match the column widths recorded for the data-entry form layout.
# (ColumnWidth is specified in "characters"; the host quantizes to a
# 6-pt-per-character grid, so these are the closest achievable values.)
$ws.Columns.Item(1).ColumnWidth = 11.6666666666667
$ws.Columns.Item(2).ColumnWidth = 10.8333333333333
$ws.Columns.Item(3).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 13.6666666666667
$ws.Columns.Item(5).ColumnWidth = 13.6666666666667
$ws.Columns.Item(6).ColumnWidth = 15.3333333333333
$ws.Columns.Item(7).ColumnWidth = 15.6666666666667

# Move the selection to A3, as left by the editing session
$ws.Range("A3").Select()
